$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Date of Meeting:  (MM/DD/YYYY)" -> "Date of Meeting:  (DD/MM/YYYY)"
# ---------------------------------------------------------------------
$dateLabelCell = $d.Tables(1).Cell(2, 1).Range
$dateLabelCell.Find.Execute("(MM/DD/YYYY)", $false, $false, $false, $false, $false, $true, 1, $false, "(DD/MM/YYYY)", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) The empty date-value cell (just a single space) gets the actual
#    meeting date appended: "16/03/2017".
# ---------------------------------------------------------------------
$dateValueCell = $d.Tables(1).Cell(2, 2)
$dateValueRange = $dateValueCell.Range
# Range.End sits right after the trailing paragraph/cell marks; back up
# one so the insertion point lands right after the existing space.
$dateInsertPoint = $d.Range($dateValueRange.End - 1, $dateValueRange.End - 1)
$dateInsertPoint.InsertAfter("16/03/2017")

# ---------------------------------------------------------------------
# 3) Fix the typo "Dr.Stphen" -> "Dr.Stephen".
# ---------------------------------------------------------------------
$facilitatorCell = $d.Tables(1).Cell(3, 2).Range
$facilitatorCell.Find.Execute("Dr.Stphen", $false, $false, $false, $false, $false, $true, 1, $false, "Dr.Stephen", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the trailing empty paragraph at
#    the end of the document to right after "Dr.Stephen" (this mirrors
#    Word's own behaviour of re-stamping _GoBack at the last edit).
#    Bookmarks.Add needs a non-empty range to actually materialise, so
#    a throw-away character is inserted, bookmarked, then removed -
#    leaving the bookmark collapsed in place exactly where Word would.
# ---------------------------------------------------------------------
$facilitatorRange = $d.Tables(1).Cell(3, 2).Range.Duplicate
$facilitatorRange.Find.Execute("Dr.Stephen") | Out-Null
$goBackAnchor = $facilitatorRange.End

$placeholder = $d.Range($goBackAnchor, $goBackAnchor)
$placeholder.InsertAfter("X")

$placeholderRange = $d.Range($goBackAnchor, $goBackAnchor + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange) | Out-Null

$placeholderRange2 = $d.Range($goBackAnchor, $goBackAnchor + 1)
$placeholderRange2.Text = ""
